# NewCustomer and CustomerList POM added
# Apply the SKU.xlsx data-entry changes to Sheet1 of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new "No of Variations" row after the "Brand Name" row (row 10) ---
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "No of Variations"
$ws.Range("B11").Value = 2

# --- Append a second variation block (rows 15-17) ---
$ws.Range("A15").Value = "variation Name"
$ws.Range("B15").Value = "test SKU automation 2"
$ws.Range("A16").Value = "Price"
$ws.Range("B16").Value = 180
$ws.Range("A17").Value = "value"
$ws.Range("B17").Value = 3

# --- Apply left alignment to the numeric "count"/"amount" cells ---
$ws.Range("B1").HorizontalAlignment = -4131
$ws.Range("B11").HorizontalAlignment = -4131
$ws.Range("B16").HorizontalAlignment = -4131
$ws.Range("B17").HorizontalAlignment = -4131

# Re-assert B1's value (keeps it numeric while the style above is applied)
$ws.Range("B1").Value = 1245

# --- Widen column B ---
$ws.Columns.Item(2).ColumnWidth = 30

# --- Reset the view: scroll to top-left and select B1 ---
$ws.Activate()
$ws.Range("B1").Select()
